# Adding final report and last comparison slide to presentation
#
# The deck currently ends with slide 12 ("Final Thoughts", sldId 334).
# We insert a brand-new "Simulation Results - Comparisons" slide right
# before it, so it becomes the new slide 12 and "Final Thoughts" slides
# down to become slide 13 (matching the sldIdLst reorder: 363 then 334).

# Helper: PowerPoint shape geometry is expressed in points over COM, but
# the OOXML stores EMU (1 pt = 12700 EMU). The interop layer marshals the
# point value through a 32-bit float before re-multiplying by 12700 and
# truncating to an integer, so a naive "emu / 12700.0" can land 1 EMU away
# from the exact target. Nudge the point value up in tiny steps until the
# round trip reproduces the exact EMU we want.
function EmuToPt($emu) {
    $base = $emu / 12700.0
    for ($k = 0; $k -lt 4000; $k++) {
        $candidate = $base + ($k * 0.0000005)
        $v32 = [float]$candidate
        $result = [int64]($v32 * 12700.0)
        if ($result -eq $emu) {
            return $candidate
        }
    }
    return $base
}

$p = $ppt.ActivePresentation

# "Title and Content" is p:sldLayoutId rId2 -> slideLayout2.xml, i.e. the
# second entry in the slide master's CustomLayouts collection (the layout
# used by nearly every other slide in this deck).
$layout = $p.SlideMaster.CustomLayouts.Item(2)

# Insert before the current last slide (index 12, "Final Thoughts") so the
# new slide becomes slide 12 and "Final Thoughts" becomes slide 13.
$newSlide = $p.Slides.AddSlide(12, $layout)

# --- Title placeholder -------------------------------------------------
$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Simulation Results - Comparisons"

# --- Body / content placeholder -----------------------------------------
$body = $newSlide.Shapes.Item(2)
$tf = $body.TextFrame

$tf.TextRange.Text = "SoC: ~12,500 samples/second"
[void]$tf.TextRange.InsertAfter("`rFPGA: ~66,000,000 samples/second")
[void]$tf.TextRange.InsertAfter("`rSimulated 100% accuracy, but implementation suffered")
[void]$tf.TextRange.InsertAfter("`rFeFET: ~416,667 samples/second")
[void]$tf.TextRange.InsertAfter("`rSimulated ~99% accuracy, but implementation suffered")
[void]$tf.TextRange.InsertAfter("`r")

# Indent levels (COM IndentLevel is 1-based; XML lvl attribute is 0-based):
#  1: SoC...                              -> lvl 0 (default)
#  2: FPGA...                             -> lvl 0 (default)
#  3: Simulated 100% accuracy...          -> lvl 1
#  4: FeFET: ~416,667 samples/second      -> lvl 0 (default)
#  5: Simulated ~99% accuracy...          -> lvl 1
#  6: (trailing empty paragraph)          -> lvl 1
$tf.TextRange.Paragraphs(3, 1).IndentLevel = 2
$tf.TextRange.Paragraphs(5, 1).IndentLevel = 2
$tf.TextRange.Paragraphs(6, 1).IndentLevel = 2

# Auto-shrink text to fit the placeholder (maps to <a:normAutofit/>).
$tf.AutoSize = 2

# Explicit position/size taken from the target geometry.
$body.Left = EmuToPt 2904565
$body.Top = EmuToPt 1825625
$body.Width = EmuToPt 8449235
$body.Height = EmuToPt 4808088
